$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: merge the 3 runs ("... and " / "AirSim" / " Simulator")
#    -- with the spell-check markers around "AirSim" -- into a single run.
#    The visible text does not change, only the run/proofErr structure, so a
#    plain Find/Replace (which no-ops on identical text) will not do the
#    merge. Instead we briefly append a sentinel character, which forces the
#    engine to rewrite the range as a single homogeneous run, then trim the
#    sentinel back off in a second pass.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleText = "Product Vision for Modelling Unmanned Aerial Swarms using Unreal Engine and AirSim Simulator"
$titleTextLen = $titleText.Length

$fullTitle = $d.Range($titleRange.Start, $titleRange.Start + $titleTextLen)
$fullTitle.Text = $titleText + "#"
$shrink = $d.Range($titleRange.Start, $titleRange.Start + $titleTextLen + 1)
$shrink.Text = $titleText

# ---------------------------------------------------------------------------
# 2) Body paragraph: replace the whole sentence with the new multi-run text.
#    We delete the paragraph's existing content (keeping the paragraph mark)
#    and re-insert it via InsertXML (flat-OPC WordprocessingML) so we get
#    exact control over run boundaries -- including keeping the leading
#    <w:tab/> as a real tab element instead of a literal "\t" character.
# ---------------------------------------------------------------------------
$bodyPara = $d.Paragraphs.Item(2)
$bodyRange = $bodyPara.Range
$clearRange = $d.Range($bodyRange.Start, $bodyRange.End - 1)
$clearRange.Delete()

$insertPoint = $d.Range($bodyRange.Start, $bodyRange.Start)

$run1 = "Our vision is to produce a simulation of a swarm of Unmanned Aerial Vehicles (UAV). A swarm is defined as a multitude of UAV moving and communicating in a unified, cohesive manner. All UAV will routinely communicate its respective position and other behavioral data to the swarm."
$run2 = " The data will then transmit to the ground control (user)."
$run3 = " The swarm will move"
$run4 = " and interact"
$run5 = " within a three-dimensional environment"
$run6 = "."
$run7 = " The swarm will gather data, such as the volume of objects in the environment. The swarm will travel along a path towards an objective"
$run8 = ". The swarm will detect objects while travelling and initiate collision avoidance if a collision is detected. This shall include avoidance of both static and moving obstacles."

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r><w:tab/><w:t>' + $run1 + '</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">' + $run2 + '</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">' + $run3 + '</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">' + $run4 + '</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">' + $run5 + '</w:t></w:r>' +
              '<w:r><w:t>' + $run6 + '</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">' + $run7 + '</w:t></w:r>' +
              '<w:r><w:t>' + $run8 + '</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertPoint.InsertXML($xml)
